# Update "想去人数" (want-to-go count) figures across the four sheets to match
# the latest scrape output (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 322
$ws.Range("F7").Value = 2115
$ws.Range("F9").Value = 47
$ws.Range("F10").Value = 1622
$ws.Range("F11").Value = 1622
$ws.Range("F12").Value = 1352
$ws.Range("F17").Value = 558
$ws.Range("F19").Value = 5
$ws.Range("F20").Value = 7194
$ws.Range("F21").Value = 7856
$ws.Range("F22").Value = 46
$ws.Range("F23").Value = 5
$ws.Range("F24").Value = 192
$ws.Range("F33").Value = 195
$ws.Range("F35").Value = 1423
$ws.Range("F36").Value = 152
$ws.Range("F37").Value = 224
$ws.Range("F40").Value = 9
$ws.Range("F41").Value = 715
$ws.Range("F42").Value = 24
$ws.Range("F44").Value = 336
$ws.Range("F48").Value = 165
$ws.Range("F49").Value = 158

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 18
$ws.Range("F5").Value = 55
$ws.Range("F11").Value = 18
$ws.Range("F17").Value = 7

# --- Sheet: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 183
$ws.Range("F5").Value = 133
$ws.Range("F6").Value = 11

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 18
$ws.Range("F5").Value = 133
$ws.Range("F7").Value = 322
$ws.Range("F10").Value = 2115
$ws.Range("F12").Value = 47
$ws.Range("F13").Value = 1622
$ws.Range("F14").Value = 1622
$ws.Range("F15").Value = 11
$ws.Range("F16").Value = 1352
$ws.Range("F20").Value = 558
$ws.Range("F22").Value = 7194
$ws.Range("F23").Value = 7856
$ws.Range("F24").Value = 46
$ws.Range("F25").Value = 5
$ws.Range("F26").Value = 192
$ws.Range("F30").Value = 195
$ws.Range("F31").Value = 1423
$ws.Range("F32").Value = 152
$ws.Range("F33").Value = 224
$ws.Range("F37").Value = 18
$ws.Range("F38").Value = 715
$ws.Range("F40").Value = 24
$ws.Range("F42").Value = 336
$ws.Range("F46").Value = 165
$ws.Range("F47").Value = 158
$ws.Range("F48").Value = 7
